# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The value 46081 (serial date 2026-02-28) is bumped to 46082 (2026-03-01)
# for every data row on the active sheet (rows 2-71).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46081) {
        $cell.Value2 = 46082
    }
}
